# Generate Report for Handback
# The 44e77bcf-... entry (row 3 in each sheet) has now been handed back,
# so its status changes from "Ready for handoff" to
# "Handed back: in sync with en-US", and the per-locale "Latest Handback
# DateTime" is stamped with the new handback timestamp.

$wb = $excel.ActiveWorkbook

$handedBack = "Handed back: in sync with en-US"

# --- Overview sheet ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $handedBack
$overview.Range("C3").Value = $handedBack

# --- zh-cn sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $handedBack
$zhcn.Range("H3").Value = "2016-03-25 09:22:15"

# --- de-de sheet ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $handedBack
$dede.Range("H3").Value = "2016-03-25 09:22:22"
